$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns keep their literal text representation
# (values like "1.007" or "0.06900" would otherwise be auto-coerced to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.703.60'
$ws.Range("E2").Value = '  -3.26%  '

$ws.Range("D3").Value = '1.872.87'
$ws.Range("E3").Value = '  -4.10%  '

$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  -0.81%  '

$ws.Range("D5").Value = '325.68'
$ws.Range("E5").Value = '  +1.41%  '

$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  -0.60%  '

$ws.Range("D7").Value = '0.4539'
$ws.Range("E7").Value = '  -4.69%  '

$ws.Range("D8").Value = '0.3893'
$ws.Range("E8").Value = '  -3.65%  '

$ws.Range("D9").Value = '48.88'
$ws.Range("E9").Value = '  -9.24%  '

$ws.Range("D10").Value = '0.08144'
$ws.Range("E10").Value = '  -3.41%  '

$ws.Range("D11").Value = '1.028'
$ws.Range("E11").Value = '  -3.01%  '

$ws.Range("D12").Value = '21.64'
$ws.Range("E12").Value = '  -2.66%  '

$ws.Range("D13").Value = '1.892.79'
$ws.Range("E13").Value = '  -3.42%  '

$ws.Range("D14").Value = '7.239'
$ws.Range("E14").Value = '  -4.79%  '

$ws.Range("D15").Value = '5.922'
$ws.Range("E15").Value = '  -4.65%  '

$ws.Range("D16").Value = '1.008'
$ws.Range("E16").Value = '  -0.76%  '

$ws.Range("D17").Value = '87.69'
$ws.Range("E17").Value = '  -1.59%  '

$ws.Range("D18").Value = '0.00001045'
$ws.Range("E18").Value = '  -2.84%  '

$ws.Range("D19").Value = '0.06572'
$ws.Range("E19").Value = '  -0.58%  '

$ws.Range("D20").Value = '17.28'
$ws.Range("E20").Value = '  -7.39%  '

$ws.Range("D21").Value = '1.004'

$ws.Range("D22").Value = '5.580'
$ws.Range("E22").Value = '  -4.26%  '

$ws.Range("D23").Value = '27.771.20'
$ws.Range("E23").Value = '  -3.08%  '

$ws.Range("D24").Value = '10.97'
$ws.Range("E24").Value = '  -4.82%  '

$ws.Range("D25").Value = '2.305'
$ws.Range("E25").Value = '  +0.35%  '

$ws.Range("D26").Value = '2.121.40'
$ws.Range("E26").Value = '  -3.10%  '

$ws.Range("D27").Value = '153.43'
$ws.Range("E27").Value = '  -0.69%  '

$ws.Range("D28").Value = '19.63'
$ws.Range("E28").Value = '  -2.74%  '

$ws.Range("D29").Value = '5.618'
$ws.Range("E29").Value = '  -5.19%  '

$ws.Range("D30").Value = '2.066'
$ws.Range("E30").Value = '  -4.10%  '

$ws.Range("D31").Value = '122.14'
$ws.Range("E31").Value = '  -1.39%  '

$ws.Range("D32").Value = '0.09472'
$ws.Range("E32").Value = '  -1.21%  '

$ws.Range("D33").Value = '0.9423'
$ws.Range("E33").Value = '  -5.75%  '

$ws.Range("D34").Value = '1.463'
$ws.Range("E34").Value = '  +2.25%  '

$ws.Range("D35").Value = '3.638'
$ws.Range("E35").Value = '  -0.91%  '

$ws.Range("D36").Value = '5.386'
$ws.Range("E36").Value = '  -3.13%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '1.246'
$ws.Range("E37").Value = '  -1.96%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02262'
$ws.Range("E38").Value = '  -3.37%  '

$ws.Range("D39").Value = '0.06041'
$ws.Range("E39").Value = '  -3.01%  '

$ws.Range("D40").Value = '8.515'
$ws.Range("E40").Value = '  -2.27%  '

$ws.Range("D41").Value = '0.6041'
$ws.Range("E41").Value = '  -2.85%  '

$ws.Range("E42").Value = '  -0.94%  '

$ws.Range("D43").Value = '10.56'
$ws.Range("E43").Value = '  -4.75%  '

$ws.Range("D44").Value = '0.1875'
$ws.Range("E44").Value = '  -2.29%  '

$ws.Range("D45").Value = '1.285'
$ws.Range("E45").Value = '  -3.46%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '12.59'
$ws.Range("E46").Value = '  -3.07%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '0.5739'
$ws.Range("E47").Value = '  -3.31%  '

$ws.Range("D48").Value = '1.968'
$ws.Range("E48").Value = '  -4.98%  '

$ws.Range("D49").Value = '3.418'
$ws.Range("E49").Value = '  +0.24%  '

$ws.Range("D50").Value = '0.06900'
$ws.Range("E50").Value = '  +1.02%  '

$ws.Range("D51").Value = '109.23'
$ws.Range("E51").Value = '  -1.74%  '
